$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.572.09"
$ws.Range("E2").Value = "  +5.02%  "

$ws.Range("D3").Value = "2.488.82"
$ws.Range("E3").Value = "  +2.65%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "322.51"
$ws.Range("E5").Value = "  +1.41%  "

$ws.Range("D6").Value = "105.53"
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  +1.78%  "

$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "0.543"
$ws.Range("E9").Value = "  +2.49%  "

$ws.Range("D10").Value = "37.96"
$ws.Range("E10").Value = "  +6.67%  "

$ws.Range("D11").Value = "0.0814"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("E12").Value = "  +1.11%  "

$ws.Range("D13").Value = "18.25"
$ws.Range("E13").Value = "  +0.77%  "

$ws.Range("D14").Value = "7.16"
$ws.Range("E14").Value = "  +1.72%  "

$ws.Range("D15").Value = "2.877.68"
$ws.Range("E15").Value = "  +2.62%  "

$ws.Range("D16").Value = "2.492.83"
$ws.Range("E16").Value = "  +2.28%  "

$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "47.434.87"
$ws.Range("E18").Value = "  +4.97%  "

$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  +3.80%  "

$ws.Range("E20").Value = "  +3.32%  "

$ws.Range("E21").Value = "  +1.82%  "

$ws.Range("D22").Value = "70.71"
$ws.Range("E22").Value = "  +2.81%  "

$ws.Range("D23").Value = "250.93"
$ws.Range("E23").Value = "  +2.80%  "

$ws.Range("D24").Value = "2.39"
$ws.Range("E24").Value = "  +5.87%  "

$ws.Range("D25").Value = "2.57"
$ws.Range("E25").Value = "  +2.71%  "

$ws.Range("D26").Value = "26.15"
$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  +4.33%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("D30").Value = "34.98"
$ws.Range("E30").Value = "  +6.37%  "

$ws.Range("D31").Value = "0.135"
$ws.Range("E31").Value = "  +6.24%  "

$ws.Range("D32").Value = "49.47"
$ws.Range("E32").Value = "  +0.80%  "

$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("D34").Value = "5.35"
$ws.Range("E34").Value = "  +2.71%  "

$ws.Range("D35").Value = "0.0781"
$ws.Range("E35").Value = "  +2.31%  "

$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("E37").Value = "  +3.57%  "

$ws.Range("E38").Value = "  +3.86%  "

$ws.Range("E39").Value = "  +4.31%  "

$ws.Range("E40").Value = "  +1.86%  "

$ws.Range("E41").Value = "  +1.70%  "

$ws.Range("D42").Value = "121.27"
$ws.Range("E42").Value = "  -3.82%  "

$ws.Range("D43").Value = "21.02"
$ws.Range("E43").Value = "  +1.85%  "

$ws.Range("E44").Value = "  +2.75%  "

$ws.Range("D45").Value = "1.962.71"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.10"
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("E48").Value = "  +1.19%  "

$ws.Range("D49").Value = "1.79"
$ws.Range("E49").Value = "  -1.74%  "

$ws.Range("D50").Value = "5.29"
$ws.Range("E50").Value = "  +12.20%  "

$ws.Range("D51").Value = "79.51"
$ws.Range("E51").Value = "  +4.08%  "
